$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect column D (numeric-looking text like "0.998", "66.406.65") from Excel
# auto-converting the assigned strings into Number-typed cells: force text format
# first, then drop back to the default "Normal" style so no stray formatting is left
# behind (keeps the produced style table identical to the original).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "66.406.65"
$ws.Range("E2").Value = "  +3.99%  "
$ws.Range("D3").Value = "3.411.18"
$ws.Range("E3").Value = "  +3.31%  "
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.41%  "
$ws.Range("D5").Value = "184.46"
$ws.Range("E5").Value = "  +3.11%  "
$ws.Range("D6").Value = "544.99"
$ws.Range("E6").Value = "  +3.47%  "
$ws.Range("D7").Value = "0.611"
$ws.Range("E7").Value = "  +1.32%  "
$ws.Range("D8").Value = "3.404.17"
$ws.Range("E8").Value = "  +3.36%  "
$ws.Range("E9").Value = "  -0.21%  "
$ws.Range("D10").Value = "0.636"
$ws.Range("E10").Value = "  +4.25%  "
$ws.Range("D11").Value = "56.05"
$ws.Range("E11").Value = "  -3.36%  "
$ws.Range("D12").Value = "0.147"
$ws.Range("E12").Value = "  +10.69%  "
$ws.Range("D13").Value = "0.0000273"
$ws.Range("E13").Value = "  +5.41%  "
$ws.Range("D14").Value = "9.37"
$ws.Range("E14").Value = "  +2.88%  "
$ws.Range("D15").Value = "3.916.95"
$ws.Range("E15").Value = "  +1.80%  "
$ws.Range("D16").Value = "0.121"
$ws.Range("E16").Value = "  +3.03%  "
$ws.Range("D17").Value = "3.379.09"
$ws.Range("E17").Value = "  +1.38%  "
$ws.Range("D18").Value = "18.22"
$ws.Range("E18").Value = "  +4.31%  "
$ws.Range("D19").Value = "66.287.26"
$ws.Range("E19").Value = "  +3.53%  "
$ws.Range("D20").Value = "11.67"
$ws.Range("E20").Value = "  +4.80%  "
$ws.Range("D21").Value = "0.996"
$ws.Range("E21").Value = "  +4.02%  "
$ws.Range("D22").Value = "407.78"
$ws.Range("E22").Value = "  +9.08%  "
$ws.Range("D23").Value = "12.02"
$ws.Range("E23").Value = "  +8.12%  "
$ws.Range("D24").Value = "4.25"
$ws.Range("E24").Value = "  +8.28%  "
$ws.Range("D25").Value = "84.23"
$ws.Range("E25").Value = "  +4.37%  "
$ws.Range("D26").Value = "3.85"
$ws.Range("E26").Value = "  +1.32%  "
$ws.Range("D27").Value = "2.88"
$ws.Range("E27").Value = "  +7.32%  "
$ws.Range("D28").Value = "6.16"
$ws.Range("E28").Value = "  +1.48%  "
$ws.Range("E29").Value = "  +1.95%  "
$ws.Range("D30").Value = "8.59"
$ws.Range("E30").Value = "  +2.77%  "
$ws.Range("D31").Value = "29.98"
$ws.Range("E31").Value = "  +3.68%  "
$ws.Range("D32").Value = "666.21"
$ws.Range("E32").Value = "  +3.02%  "
$ws.Range("D33").Value = "6.85"
$ws.Range("E33").Value = "  +3.18%  "
$ws.Range("D34").Value = "11.57"
$ws.Range("E34").Value = "  +2.68%  "
$ws.Range("D35").Value = "0.110"
$ws.Range("E35").Value = "  +3.50%  "
$ws.Range("D36").Value = "58.39"
$ws.Range("E36").Value = "  -1.41%  "
$ws.Range("B37").Value = "InjectiveProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D37").Value = "38.40"
$ws.Range("E37").Value = "  +4.73%  "
$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").Value = "0.0₃0814"
$ws.Range("E38").Value = "  +16.15%  "
$ws.Range("B39").Value = "TheGraph"
$ws.Range("C39").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D39").Value = "0.400"
$ws.Range("E39").Value = "  +2.64%  "
$ws.Range("B40").Value = "Dai"
$ws.Range("C40").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("D41").Value = "2.81"
$ws.Range("E41").Value = "  +13.52%  "
$ws.Range("D42").Value = "0.133"
$ws.Range("E42").Value = "  +6.15%  "
$ws.Range("D43").Value = "3.31"
$ws.Range("E43").Value = "  +18.68%  "
$ws.Range("D44").Value = "0.996"
$ws.Range("E44").Value = "  -0.55%  "
$ws.Range("D45").Value = "3.027.31"
$ws.Range("E45").Value = "  +2.60%  "
$ws.Range("D46").Value = "2.86"
$ws.Range("E46").Value = "  +5.76%  "
$ws.Range("D47").Value = "0.0417"
$ws.Range("E47").Value = "  +4.81%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").Value = "3.24"
$ws.Range("E48").Value = "  +4.89%  "
$ws.Range("B49").Value = "WEMIXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").Value = "2.75"
$ws.Range("E49").Value = "  +3.89%  "
$ws.Range("D50").Value = "8.82"
$ws.Range("E50").Value = "  +11.81%  "
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").Value = "2.64"
$ws.Range("E51").Value = "  +6.31%  "

$ws.Range("D2:D51").Style = "Normal"

